# "fixed email app, now reads from excel file"
# Append the next email address as a new row, turning it into a live
# mailto: hyperlink (Excel's built-in "Hyperlink" cell style), then move
# the selection down to the next empty row ready for the following entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newEmail = "triminhtran2797@gmail.com"

$cell = $ws.Range("A3")
$cell.Value = $newEmail
$ws.Hyperlinks.Add($cell, "mailto:$newEmail")

$ws.Range("A4").Select()
